$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

    $ws.Range("B2").Value = 15.53426190971729
    $ws.Range("C2").Value = 12.42272858807831
    $ws.Range("D2").Value = 5.968553232837087
    $ws.Range("E2").Value = 16.52896572815041
    $ws.Range("G2").Value = 34.82032796537018
    $ws.Range("H2").Value = 14.64808963408755
    $ws.Range("I2").Value = 20.054024083525
    $ws.Range("B3").Value = 14.78388018740584
    $ws.Range("C3").Value = 11.60122308761781
    $ws.Range("D3").Value = 5.848011660073984
    $ws.Range("E3").Value = 15.5837134517504
    $ws.Range("G3").Value = 34.02863113231896
    $ws.Range("H3").Value = 14.62321074893769
    $ws.Range("I3").Value = 20.0672728556703
    $ws.Range("B4").Value = 14.3061591948401
    $ws.Range("C4").Value = 11.09958326520868
    $ws.Range("D4").Value = 5.774746439216896
    $ws.Range("E4").Value = 14.97905080547837
    $ws.Range("G4").Value = 33.55240806141679
    $ws.Range("H4").Value = 14.61308305399272
    $ws.Range("I4").Value = 20.08349520147614
    $ws.Range("B5").Value = 14.10748008055997
    $ws.Range("C5").Value = 10.88880438345576
    $ws.Range("D5").Value = 5.745122979709981
    $ws.Range("E5").Value = 14.72681139110883
    $ws.Range("G5").Value = 33.36115491904042
    $ws.Range("H5").Value = 14.6102455655456
    $ws.Range("I5").Value = 20.09212050972133
    $ws.Range("B6").Value = 14.07425634625698
    $ws.Range("C6").Value = 10.85342428068708
    $ws.Range("D6").Value = 5.740219391171711
    $ws.Range("E6").Value = 14.68458376907187
    $ws.Range("G6").Value = 33.32957677312164
    $ws.Range("H6").Value = 14.60985212099266
    $ws.Range("I6").Value = 20.0936738194254
    $ws.Range("B7").Value = 14.30349556782921
    $ws.Range("C7").Value = 11.09676621172561
    $ws.Range("D7").Value = 5.774345925453486
    $ws.Range("E7").Value = 14.97567224501798
    $ws.Range("G7").Value = 33.5498169411624
    $ws.Range("H7").Value = 14.61303957249589
    $ws.Range("I7").Value = 20.08360339490188
    $ws.Range("B8").Value = 15.27920284607282
    $ws.Range("C8").Value = 12.14293145523802
    $ws.Range("D8").Value = 5.92686200132132
    $ws.Range("E8").Value = 16.20822660877739
    $ws.Range("G8").Value = 34.54549659249503
    $ws.Range("H8").Value = 14.6384404281166
    $ws.Range("I8").Value = 20.05690270928379
    $ws.Range("B9").Value = 17.04804094089212
    $ws.Range("C9").Value = 14.04090658862258
    $ws.Range("D9").Value = 6.229940523580388
    $ws.Range("E9").Value = 18.52423419746841
    $ws.Range("G9").Value = 36.56205376678879
    $ws.Range("H9").Value = 14.7292479132541
    $ws.Range("I9").Value = 20.06948513428042
    $ws.Range("B10").Value = 18.24880853026226
    $ws.Range("C10").Value = 15.28468413404332
    $ws.Range("D10").Value = 6.452520187193766
    $ws.Range("E10").Value = 20.18903791592331
    $ws.Range("G10").Value = 38.06373092124111
    $ws.Range("H10").Value = 14.82106027004779
    $ws.Range("I10").Value = 20.11924866805255
    $ws.Range("B11").Value = 18.77193702301494
    $ws.Range("C11").Value = 15.8181564507246
    $ws.Range("D11").Value = 6.553251061043894
    $ws.Range("E11").Value = 20.90472012969369
    $ws.Range("G11").Value = 38.74743353645941
    $ws.Range("H11").Value = 14.86826687020556
    $ws.Range("I11").Value = 20.15084282183847
    $ws.Range("B12").Value = 18.96660074235519
    $ws.Range("C12").Value = 16.01554799009641
    $ws.Range("D12").Value = 6.591280558529864
    $ws.Range("E12").Value = 21.16978312090025
    $ws.Range("G12").Value = 39.00611897728928
    $ws.Range("H12").Value = 14.88692164183589
    $ws.Range("I12").Value = 20.16410416687249
    $ws.Range("B13").Value = 18.92483069932225
    $ws.Range("C13").Value = 15.97324138090954
    $ws.Range("D13").Value = 6.583095959503771
    $ws.Range("E13").Value = 21.11296109862336
    $ws.Range("G13").Value = 38.95042075011318
    $ws.Range("H13").Value = 14.88286944746616
    $ws.Range("I13").Value = 20.16119027374942
    $ws.Range("B14").Value = 18.78802150591827
    $ws.Range("C14").Value = 15.83448847206953
    $ws.Range("D14").Value = 6.556382299151796
    $ws.Range("E14").Value = 20.92664605647586
    $ws.Range("G14").Value = 38.76872177717821
    $ws.Range("H14").Value = 14.8697860308115
    $ws.Range("I14").Value = 20.15190779435038
    $ws.Range("B15").Value = 18.70377180317484
    $ws.Range("C15").Value = 15.74889704547596
    $ws.Range("D15").Value = 6.54000326972706
    $ws.Range("E15").Value = 20.81174893700022
    $ws.Range("G15").Value = 38.65738887222074
    $ws.Range("H15").Value = 14.8618733278702
    $ws.Range("I15").Value = 20.14639119906327
    $ws.Range("B16").Value = 18.21414593846746
    $ws.Range("C16").Value = 15.24917215312996
    $ws.Range("D16").Value = 6.44592321397618
    $ws.Range("E16").Value = 20.14143165930408
    $ws.Range("G16").Value = 38.01903489490483
    $ws.Range("H16").Value = 14.81808440688571
    $ws.Range("I16").Value = 20.11736499590156
    $ws.Range("B17").Value = 17.90777266236588
    $ws.Range("C17").Value = 14.93434553237988
    $ws.Range("D17").Value = 6.3880453243446
    $ws.Range("E17").Value = 19.71957458409835
    $ws.Range("G17").Value = 37.62735463482127
    $ws.Range("H17").Value = 14.79261267127431
    $ws.Range("I17").Value = 20.10185970430567
    $ws.Range("B18").Value = 17.72938518023917
    $ws.Range("C18").Value = 14.75022114628256
    $ws.Range("D18").Value = 6.354708650798845
    $ws.Range("E18").Value = 19.47301141714143
    $ws.Range("G18").Value = 37.4021423431797
    $ws.Range("H18").Value = 14.77847442934638
    $ws.Range("I18").Value = 20.09378407209918
    $ws.Range("B19").Value = 17.66861731532074
    $ws.Range("C19").Value = 14.6873561161494
    $ws.Range("D19").Value = 6.34341463717925
    $ws.Range("E19").Value = 19.38885456889157
    $ws.Range("G19").Value = 37.32591113258987
    $ws.Range("H19").Value = 14.77377556072826
    $ws.Range("I19").Value = 20.09119419722452
    $ws.Range("B20").Value = 17.94061204157287
    $ws.Range("C20").Value = 14.96817426348869
    $ws.Range("D20").Value = 6.394211669700151
    $ws.Range("E20").Value = 19.76488763860476
    $ws.Range("G20").Value = 37.66904447117383
    $ws.Range("H20").Value = 14.79527116875002
    $ws.Range("I20").Value = 20.10342298804605
    $ws.Range("B21").Value = 18.82829965795709
    $ws.Range("C21").Value = 15.87536882145062
    $ws.Range("D21").Value = 6.564232180955802
    $ws.Range("E21").Value = 20.98153248702633
    $ws.Range("G21").Value = 38.82209940258156
    $ws.Range("H21").Value = 14.87360785212799
    $ws.Range("I21").Value = 20.15459901052074
    $ws.Range("B22").Value = 19.38840257596116
    $ws.Range("C22").Value = 16.44133937044438
    $ws.Range("D22").Value = 6.674662525573936
    $ws.Range("E22").Value = 21.74201479648275
    $ws.Range("G22").Value = 39.5742898544899
    $ws.Range("H22").Value = 14.92934120714339
    $ws.Range("I22").Value = 20.19561010465687
    $ws.Range("B23").Value = 19.0913315656208
    $ws.Range("C23").Value = 16.14172611609272
    $ws.Range("D23").Value = 6.615799378781333
    $ws.Range("E23").Value = 21.33929036961712
    $ws.Range("G23").Value = 39.17305468610778
    $ws.Range("H23").Value = 14.89918186305172
    $ws.Range("I23").Value = 20.17302701382899
    $ws.Range("B24").Value = 17.9257723560033
    $ws.Range("C24").Value = 14.95289003405701
    $ws.Range("D24").Value = 6.391424054716765
    $ws.Range("E24").Value = 19.74441416807204
    $ws.Range("G24").Value = 37.65019656681444
    $ws.Range("H24").Value = 14.79406768633027
    $ws.Range("I24").Value = 20.10271361625842
    $ws.Range("B25").Value = 16.58621826143823
    $ws.Range("C25").Value = 13.55403441873346
    $ws.Range("D25").Value = 6.147788141916061
    $ws.Range("E25").Value = 17.8739712882458
    $ws.Range("G25").Value = 36.01165168786137
    $ws.Range("H25").Value = 14.70027621358755
    $ws.Range("I25").Value = 20.05902517514695
